$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C6").Value = "t"
$ws.Range("C7").Value = "t"
$ws.Range("C10").Value = "t"
